$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End | Animal Glue
$ws.Range("H5").Value = 229.88889
$ws.Range("I5").Value = 233.625
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 233.625
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -118.625
$ws.Range("N5").Value = -430

# Row 18: You Grow, Girl | Growth Formula Beta
$ws.Range("H18").Value = 955
$ws.Range("I18").Value = 955
$ws.Range("K18").Value = 955
$ws.Range("M18").Value = -671

# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Range("H43").Value = 924.9286
$ws.Range("I43").Value = 421.2857
$ws.Range("J43").Value = 1428.5714
$ws.Range("K43").Value = 421.2857
$ws.Range("L43").Value = 1428.5714
$ws.Range("M43").Value = -352.2857
$ws.Range("N43").Value = -1566.5714

# Row 111: An Eye for Healing | Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 2900
$ws.Range("I111").Value = 2900
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 8700
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -5633

# Row 115: 5-bell Energy | Competent Craftsman's Syrup
$ws.Range("H115").Value = 492.5
$ws.Range("I115").Value = 492.5
$ws.Range("K115").Value = 1477.5
$ws.Range("M115").Value = 89.5

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 905.1836499999999
$ws.Range("J129").Value = 923.5111000000001
$ws.Range("L129").Value = 2770.5333
$ws.Range("N129").Value = -12770.5333

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth | Bronze Rivets
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -88

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 12255.372
$ws.Range("I32").Value = 13043.191
$ws.Range("J32").Value = 2998.5
$ws.Range("K32").Value = 13043.191
$ws.Range("L32").Value = 2998.5
$ws.Range("M32").Value = -12756.191
$ws.Range("N32").Value = -3572.5

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 2258.2083
$ws.Range("I61").Value = 1705.3889
$ws.Range("J61").Value = 3916.6667
$ws.Range("K61").Value = 1705.3889
$ws.Range("L61").Value = 3916.6667
$ws.Range("M61").Value = -1493.3889
$ws.Range("N61").Value = -4340.6667

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1197.7778
$ws.Range("I74").Value = 836
$ws.Range("K74").Value = 836
$ws.Range("M74").Value = 38

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1197.7778
$ws.Range("I77").Value = 836
$ws.Range("K77").Value = 4180
$ws.Range("M77").Value = 188

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 1961.4736
$ws.Range("I122").Value = 1816.9333
$ws.Range("K122").Value = 5450.7999
$ws.Range("M122").Value = -3000.7999

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2258.2083
$ws.Range("I136").Value = 1705.3889
$ws.Range("J136").Value = 3916.6667
$ws.Range("K136").Value = 5116.1667
$ws.Range("L136").Value = 11750.0001
$ws.Range("M136").Value = -2566.1667
$ws.Range("N136").Value = -16850.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences | Bronze Rivets
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -85

# Row 22: Riveting Run | Iron Rivets
$ws.Range("H22").Value = 6151.1113
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 59: Pop That Top | Cobalt Raising Hammer
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").ClearContents()
$ws.Range("N59").Value = 0

$ws = $wb.Worksheets.Item("CRP")
# Row 8: Bows for the Boys | Maple Longbow
$ws.Range("H8").Value = 15954.5
$ws.Range("I8").Value = 1009
$ws.Range("J8").Value = 30900
$ws.Range("K8").Value = 1009
$ws.Range("L8").Value = 30900
$ws.Range("M8").Value = -869
$ws.Range("N8").Value = -31180

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2027.9767
$ws.Range("I31").Value = 1308.3438
$ws.Range("J31").Value = 4121.4546
$ws.Range("K31").Value = 1308.3438
$ws.Range("L31").Value = 4121.4546
$ws.Range("M31").Value = -1013.3438
$ws.Range("N31").Value = -4711.4546

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2027.9767
$ws.Range("I34").Value = 1308.3438
$ws.Range("J34").Value = 4121.4546
$ws.Range("K34").Value = 1308.3438
$ws.Range("L34").Value = 4121.4546
$ws.Range("M34").Value = -1106.3438
$ws.Range("N34").Value = -4525.4546

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2191.6
$ws.Range("I58").Value = 2492
$ws.Range("J58").Value = 990
$ws.Range("K58").Value = 2492
$ws.Range("L58").Value = 990
$ws.Range("M58").Value = -2289
$ws.Range("N58").Value = -1396

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 2532.2
$ws.Range("I99").Value = 2448.4707
$ws.Range("K99").Value = 2448.4707
$ws.Range("M99").Value = -950.4706999999999

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 2333
$ws.Range("I122").Value = 2148.2173
$ws.Range("J122").Value = 2864.25
$ws.Range("K122").Value = 6444.651899999999
$ws.Range("L122").Value = 8592.75
$ws.Range("M122").Value = -3994.651899999999
$ws.Range("N122").Value = -13492.75

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 2532.2
$ws.Range("I126").Value = 2448.4707
$ws.Range("K126").Value = 7345.4121
$ws.Range("M126").Value = -4875.4121

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 2715.3103
$ws.Range("I132").Value = 1837.3
$ws.Range("J132").Value = 4666.4443
$ws.Range("K132").Value = 5511.9
$ws.Range("L132").Value = 13999.3329
$ws.Range("M132").Value = -2981.9
$ws.Range("N132").Value = -19059.3329

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2074.5312
$ws.Range("I134").Value = 1922.5
$ws.Range("J134").Value = 2409
$ws.Range("K134").Value = 5767.5
$ws.Range("L134").Value = 7227
$ws.Range("M134").Value = -3232.5
$ws.Range("N134").Value = -12297

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2191.6
$ws.Range("I136").Value = 2492
$ws.Range("J136").Value = 990
$ws.Range("K136").Value = 7476
$ws.Range("L136").Value = 2970
$ws.Range("M136").Value = -4926
$ws.Range("N136").Value = -8070

$ws = $wb.Worksheets.Item("CUL")
# Row 25: Flakes for Friends | Apple Tart
$ws.Range("H25").Value = 1474.1428
$ws.Range("J25").Value = 4799.5
$ws.Range("L25").Value = 14398.5
$ws.Range("N25").Value = -14736.5

# Row 30: Picnic Panic | Apple Tart
$ws.Range("H30").Value = 1474.1428
$ws.Range("J30").Value = 4799.5
$ws.Range("L30").Value = 14398.5
$ws.Range("N30").Value = -14602.5

# Row 99: A Shorlonging for the Familiar | Shorlog
$ws.Range("H99").Value = 8166.6665
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 8166.6665
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").Value = 24499.9995
$ws.Range("N99").Value = -28991.9995

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 2846.7058
$ws.Range("I122").Value = 1598.4445
$ws.Range("J122").Value = 4251
$ws.Range("K122").Value = 4795.333500000001
$ws.Range("L122").Value = 12753
$ws.Range("M122").Value = -2345.333500000001
$ws.Range("N122").Value = -17653

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 3709.9092
$ws.Range("I7").Value = 4160.8
$ws.Range("J7").Value = 3334.1667
$ws.Range("K7").Value = 4160.8
$ws.Range("L7").Value = 3334.1667
$ws.Range("M7").Value = -4048.8
$ws.Range("N7").Value = -3558.1667

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 5039.615
$ws.Range("I40").Value = 10900
$ws.Range("J40").Value = 3281.5
$ws.Range("K40").Value = 10900
$ws.Range("L40").Value = 3281.5
$ws.Range("M40").Value = -10764
$ws.Range("N40").Value = -3553.5

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 10531853
$ws.Range("I122").Value = 5246.231
$ws.Range("K122").Value = 15738.693
$ws.Range("M122").Value = -13288.693

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 3709.9092
$ws.Range("I126").Value = 4160.8
$ws.Range("J126").Value = 3334.1667
$ws.Range("K126").Value = 12482.4
$ws.Range("L126").Value = 10002.5001
$ws.Range("M126").Value = -10012.4
$ws.Range("N126").Value = -14942.5001

$ws = $wb.Worksheets.Item("WVR")
# Row 74: Clothing the Naked Truth | Ramie Robe of Casting
$ws.Range("H74").Value = 6054.4443
$ws.Range("I74").Value = 2569
$ws.Range("J74").Value = 6490.125
$ws.Range("K74").Value = 2569
$ws.Range("L74").Value = 6490.125
$ws.Range("M74").Value = -1633
$ws.Range("N74").Value = -8362.125

# Row 77: When in Robes (L) | Ramie Robe of Casting
$ws.Range("H77").Value = 6054.4443
$ws.Range("I77").Value = 2569
$ws.Range("J77").Value = 6490.125
$ws.Range("K77").Value = 7707
$ws.Range("L77").Value = 19470.375
$ws.Range("M77").Value = -3027
$ws.Range("N77").Value = -28830.375
